$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1726.5186
$ws.Range("I15").Value = 1726.5186
$ws.Range("K15").Value = 5179.5558
$ws.Range("M15").Value = -5010.5558

# Row 33
$ws.Range("H33").Value = 210.4
$ws.Range("J33").Value = 225
$ws.Range("L33").Value = 225
$ws.Range("N33").Value = -683

# Row 86
$ws.Range("H86").Value = 13713.25
$ws.Range("I86").Value = 21374.5
$ws.Range("J86").Value = 6052
$ws.Range("K86").Value = 21374.5
$ws.Range("L86").Value = 6052
$ws.Range("M86").Value = -20251.5
$ws.Range("N86").Value = -8298

# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 89
$ws.Range("H89").Value = 13713.25
$ws.Range("I89").Value = 21374.5
$ws.Range("J89").Value = 6052
$ws.Range("K89").Value = 106872.5
$ws.Range("L89").Value = 30260
$ws.Range("M89").Value = -101256.5
$ws.Range("N89").Value = -41492

# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 92
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()

# Row 132
$ws.Range("H132").Value = 3218.2632
$ws.Range("I132").Value = 3230.389
$ws.Range("K132").Value = 9691.167000000001
$ws.Range("M132").Value = -7161.167000000001

# Row 138
$ws.Range("H138").Value = 2316.348
$ws.Range("I138").Value = 1653.9231
$ws.Range("J138").Value = 2577.303
$ws.Range("K138").Value = 4961.7693
$ws.Range("L138").Value = 7731.909
$ws.Range("M138").Value = 178.2307000000001
$ws.Range("N138").Value = -18011.909


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 5981.7856
$ws.Range("I45").Value = 7131.7
$ws.Range("K45").Value = 7131.7
$ws.Range("M45").Value = -6754.7

# Row 61
$ws.Range("H61").Value = 125002584
$ws.Range("I61").Value = 142858820
$ws.Range("K61").Value = 142858820
$ws.Range("M61").Value = -142858608

# Row 102
$ws.Range("H102").Value = 5264519.5
$ws.Range("J102").Value = 2388.3333
$ws.Range("L102").Value = 2388.3333
$ws.Range("N102").Value = -5632.3333

# Row 122
$ws.Range("H122").Value = 4741.8623
$ws.Range("I122").Value = 3960.56
$ws.Range("J122").Value = 9625
$ws.Range("K122").Value = 11881.68
$ws.Range("L122").Value = 28875
$ws.Range("M122").Value = -9431.68
$ws.Range("N122").Value = -33775

# Row 132
$ws.Range("H132").Value = 3706477
$ws.Range("I132").Value = 3848983.8
$ws.Range("K132").Value = 11546951.4
$ws.Range("M132").Value = -11544421.4

# Row 136
$ws.Range("H136").Value = 125002584
$ws.Range("I136").Value = 142858820
$ws.Range("K136").Value = 428576460
$ws.Range("M136").Value = -428573910


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 413.75
$ws.Range("I5").Value = 413.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 413.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -300.75
$ws.Range("N5").ClearContents()

# Row 117
$ws.Range("H117").Value = 26992
$ws.Range("J117").Value = 26992
$ws.Range("L117").Value = 26992
$ws.Range("N117").Value = -36170


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1098268.2
$ws.Range("I16").Value = 1567765.2
$ws.Range("J16").Value = 2775
$ws.Range("K16").Value = 1567765.2
$ws.Range("L16").Value = 2775
$ws.Range("M16").Value = -1567478.2
$ws.Range("N16").Value = -3349

# Row 31
$ws.Range("H31").Value = 6900
$ws.Range("I31").Value = 8000
$ws.Range("K31").Value = 8000
$ws.Range("M31").Value = -7705

# Row 34
$ws.Range("H34").Value = 6900
$ws.Range("I34").Value = 8000
$ws.Range("K34").Value = 8000
$ws.Range("M34").Value = -7798

# Row 107
$ws.Range("H107").Value = 638379.4
$ws.Range("I107").Value = 1235214.6
$ws.Range("J107").Value = 101227.6
$ws.Range("K107").Value = 1235214.6
$ws.Range("L107").Value = 101227.6
$ws.Range("M107").Value = -1233294.6
$ws.Range("N107").Value = -105067.6

# Row 113
$ws.Range("H113").Value = 1098268.2
$ws.Range("I113").Value = 1567765.2
$ws.Range("J113").Value = 2775
$ws.Range("K113").Value = 1567765.2
$ws.Range("L113").Value = 2775
$ws.Range("M113").Value = -1565595.2
$ws.Range("N113").Value = -7115


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 199.11111
$ws.Range("J17").Value = 350.5
$ws.Range("L17").Value = 1051.5
$ws.Range("N17").Value = -1389.5

# Row 38
$ws.Range("H38").Value = 189.92308
$ws.Range("I38").Value = 294
$ws.Range("J38").Value = 124.875
$ws.Range("K38").Value = 882
$ws.Range("L38").Value = 374.625
$ws.Range("M38").Value = -535
$ws.Range("N38").Value = -1068.625

# Row 113
$ws.Range("H113").Value = 167600
$ws.Range("I113").Value = 250649.25
$ws.Range("J113").Value = 1501.5
$ws.Range("K113").Value = 751947.75
$ws.Range("L113").Value = 4504.5
$ws.Range("M113").Value = -749777.75
$ws.Range("N113").Value = -8844.5


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2282.6667
$ws.Range("I80").Value = 2282.6667
$ws.Range("K80").Value = 2282.6667
$ws.Range("M80").Value = -1284.6667

# Row 83
$ws.Range("H83").Value = 2282.6667
$ws.Range("I83").Value = 2282.6667
$ws.Range("K83").Value = 11413.3335
$ws.Range("M83").Value = -6421.333500000001

# Row 132
$ws.Range("H132").Value = 4467253.5
$ws.Range("I132").Value = 4810735.5
$ws.Range("K132").Value = 14432206.5
$ws.Range("M132").Value = -14429676.5


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1602.85
$ws.Range("I16").Value = 1188.8667
$ws.Range("K16").Value = 1188.8667
$ws.Range("M16").Value = -1018.8667

# Row 111
$ws.Range("H111").Value = 74380
$ws.Range("J111").Value = 74380
$ws.Range("L111").Value = 74380
$ws.Range("N111").Value = -82560


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 39495.75
$ws.Range("J41").Value = 39495.75
$ws.Range("L41").Value = 39495.75
$ws.Range("N41").Value = -40275.75

# Row 116
$ws.Range("H116").Value = 82249
$ws.Range("J116").Value = 82249
$ws.Range("L116").Value = 82249
$ws.Range("N116").Value = -91427

# Row 132
$ws.Range("H132").Value = 15630586
$ws.Range("I132").Value = 22729658
$ws.Range("K132").Value = 68188974
$ws.Range("M132").Value = -68186444

